# 13_16_stimuli.xlsx -- "more work towards final product"
#
# Fills in the still-blank kind/carrier (and pair_kind) columns for the
# practice rows, the four "generic" word rows, and the 8 new
# unique_video/unique_audio rows that were added lower in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Practice rows (p1..p4, rows 2-5): carrier column D gets the matching
#    pair word (can/where/do/look), mirroring column K on the same row.
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# -- Generic word rows (1-4, sheet rows 6-9): pair_kind column J gets the
#    unique_video/unique_audio label for that pair.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# -- New numbers 9-16 (sheet rows 14-21): kind (C) + carrier (D), following
#    the same unique_video/unique_audio + can/where/do/look pattern.
$rows = @(
    @{ Row = 14; Kind = "unique_video"; Carrier = "look" },
    @{ Row = 15; Kind = "unique_video"; Carrier = "look" },
    @{ Row = 16; Kind = "unique_video"; Carrier = "where" },
    @{ Row = 17; Kind = "unique_video"; Carrier = "where" },
    @{ Row = 18; Kind = "unique_audio"; Carrier = "can" },
    @{ Row = 19; Kind = "unique_audio"; Carrier = "can" },
    @{ Row = 20; Kind = "unique_audio"; Carrier = "do" },
    @{ Row = 21; Kind = "unique_audio"; Carrier = "do" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.Kind      # column C
    $ws.Cells.Item($r.Row, 4).Value = $r.Carrier   # column D
}
